{"js": "// Update the date paragraph and the 25 practice-problem cells in the table.\nconst body = context.document.body;\n\n// --- 1) Date heading paragraph: \"2024-04-14 Sunday\" -> \"2024-04-15 Monday\" ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.load(\"text\");\nawait context.sync();\n\nif (dateParagraph.text.indexOf(\"2024-04-14 Sunday\") !== -1) {\n  dateParagraph.getRange().insertText(\"2024-04-15 Monday\", \"Replace\");\n}\n\n// --- 2) Table cells: replace each answer with its new value, by position ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Only rows 0, 4, 8, 12, 16 hold problems; the rest are blank spacer rows.\n// Values are listed row-by-row, left to right (col 0..4), matching the\n// document's reading order.\nconst newValues = {\n  0: [\"93\u00f73=31, 0\", \"34\u00f72=17, 0\", \"16\u00f75=3, 1\", \"96\u00f73=32, 0\", \"73\u00f76=12, 1\"],\n  4: [\"64\u00f78=8, 0\", \"33\u00f77=4, 5\", \"38\u00f79=4, 2\", \"63\u00f79=7, 0\", \"88\u00f72=44, 0\"],\n  8: [\"29\u00f75=5, 4\", \"72\u00f79=8, 0\", \"63\u00f76=10, 3\", \"91\u00f73=30, 1\", \"63\u00f75=12, 3\"],\n  12: [\"13\u00f74=3, 1\", \"31\u00f72=15, 1\", \"88\u00f78=11, 0\", \"76\u00f73=25, 1\", \"74\u00f79=8, 2\"],\n  16: [\"26\u00f78=3, 2\", \"65\u00f75=13, 0\", \"59\u00f75=11, 4\", \"66\u00f72=33, 0\", \"51\u00f78=6, 3\"],\n};\n\nfor (const rowIndexStr of Object.keys(newValues)) {\n  const rowIndex = Number(rowIndexStr);\n  const rowValues = newValues[rowIndex];\n  for (let col = 0; col < rowValues.length; col++) {\n    table.getCell(rowIndex, col).value = rowValues[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the 25 practice-problem answers in the table.\n$d = $word.ActiveDocument\n\n# --- 1) Date heading paragraph: \"2024-04-14 Sunday\" -> \"2024-04-15 Monday\" ---\n$dateParagraph = $d.Paragraphs.Item(1)\nif ($dateParagraph.Range.Text -like \"*2024-04-14 Sunday*\") {\n    $dateParagraph.Range.Text = \"2024-04-15 Monday\"\n}\n\n# --- 2) Table cells: replace each answer with its new value, by position ---\n# Only rows 1, 5, 9, 13, 17 (1-indexed) hold problems; the rest are blank\n# spacer rows. Values are listed row-by-row, left to right (col 1..5).\n$t = $d.Tables.Item(1)\n\n$newValues = @{\n    1  = @(\"93\u00f73=31, 0\", \"34\u00f72=17, 0\", \"16\u00f75=3, 1\", \"96\u00f73=32, 0\", \"73\u00f76=12, 1\")\n    5  = @(\"64\u00f78=8, 0\", \"33\u00f77=4, 5\", \"38\u00f79=4, 2\", \"63\u00f79=7, 0\", \"88\u00f72=44, 0\")\n    9  = @(\"29\u00f75=5, 4\", \"72\u00f79=8, 0\", \"63\u00f76=10, 3\", \"91\u00f73=30, 1\", \"63\u00f75=12, 3\")\n    13 = @(\"13\u00f74=3, 1\", \"31\u00f72=15, 1\", \"88\u00f78=11, 0\", \"76\u00f73=25, 1\", \"74\u00f79=8, 2\")\n    17 = @(\"26\u00f78=3, 2\", \"65\u00f75=13, 0\", \"59\u00f75=11, 4\", \"66\u00f72=33, 0\", \"51\u00f78=6, 3\")\n}\n\nforeach ($rowIndex in $newValues.Keys) {\n    $rowValues = $newValues[$rowIndex]\n    for ($col = 1; $col -le $rowValues.Length; $col++) {\n        $t.Cell($rowIndex, $col).Range.Text = $rowValues[$col - 1]\n    }\n}\n"}
